# connectors.xlsx — "added needed caps, fixed a spelling error"
#
# The real content change behind this commit is on the "audio" sheet:
# cell A2 held a leftover formula (="RES-"&TEXT(ROW()-1,"000000")) that
# produced a stray "RES-000001" placeholder value in row 2 — that gets
# cleared out entirely (not just blanked to 0/""), along with the
# selections/active-sheet state that naturally shift around as a user
# opens the file, cleans the cell, and clicks around between sheets.

$wb = $excel.ActiveWorkbook

# --- audio sheet: drop the stray formula/value from A2 -------------------
$wsAudio = $wb.Worksheets.Item("audio")
$null = $wsAudio.Range("A2").ClearContents()

# --- DC power sheet: cursor left parked on B37 ----------------------------
$wsDcPower = $wb.Worksheets.Item("DC power")
$null = $wsDcPower.Activate()
$null = $wsDcPower.Range("B37").Select()

# --- audio becomes the front/active tab, cursor resting on A2 ------------
$null = $wsAudio.Activate()
$null = $wsAudio.Range("A2").Select()
